$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 'សំណើម៖ វាគឺជាសមាមាត្រនៃសម្ពាធផ្នែកនៃចំហាយទឹកទៅនឹងសម្ពាធចំហាយទឹកដែលមានលំនឹងនៅសីតុណ្ហភាពដែលបានផ្តល់ឱ្យ។ ដើម្បីសម្រេចបាននូវសីតុណ្ហភាពល្អបំផុតនៅក្នុងផ្ទះ វាជាការសំខាន់ណាស់ដែលសំណើមដែលទាក់ទងនៃខ្យល់គឺលើសពី 40% ហើយវាមិនលើសពី 60% ទេ។ សមាមាត្រខ្ពស់នៃសំណើមក៏ជួយដល់ការវិវត្តនៃរោគផ្លូវដង្ហើម និងការបង្កើតផ្សិតផងដែរ។'
$ws.Range("D4").Value = 'CO2៖ កាបូនឌីអុកស៊ីតគឺជាឧស្ម័នគ្មានក្លិន គ្មានពណ៌ គ្មានជាតិអាស៊ីតបន្តិច និងមិនងាយឆេះ។ ការស្រូបចូលដែលមានកំហាប់ខ្ពស់អាចបណ្តាលឱ្យមានខ្យល់ចេញចូលខ្លាំង បាត់បង់ស្មារតី tachycardia និងឈឺក្បាល។ ប្រសិនបើ​ការ​ប៉ះពាល់​យូរ​ទៅ វា​អាច​បង្ក​ឱ្យ​មាន​ការ​ប្រែប្រួល​ក្នុង​ការ​រំលាយ​អាហារ​របស់​មនុស្ស​។'
$ws.Range("E4").Value = 'CO: កាបូនម៉ូណូអុកស៊ីតគឺជាឧស្ម័នគ្មានក្លិន គ្មានពណ៌ គ្មានរសជាតិ គ្មានជាតិពុល និងងាយឆេះ បើទោះបីជាវាមិនមែនជាសារធាតុដែលឆាប់ខឹងក៏ដោយ ដូច្នេះការប៉ះពាល់អាចត្រូវបានគេកត់សម្គាល់ទាំងស្រុង។ អាចនាំអោយមានការច្របូកច្របល់ផ្លូវចិត្ត វិលមុខ ឈឺក្បាល ចង្អោរ ខ្សោយ និងបាត់បង់ស្មារតី។ ប្រសិនបើការប៉ះពាល់យូរកើតឡើង ប្រព័ន្ធប្រសាទ និងសរសៃឈាមបេះដូងអាចនឹងរងផលប៉ះពាល់ ដែលនាំឱ្យមានបញ្ហាប្រព័ន្ធប្រសាទ និងបេះដូង។'
$ws.Range("F4").Value = 'PM 1.0: ពួកវាជាភាគល្អិតព្យួរដែលមានអង្កត់ផ្ចិតខ្យល់រហូតដល់ 0.1 µm ដែលហៅថាភាគល្អិតជ្រុល ឬប្រភាគ ultrafine ។ ពួកវាមាននៅក្នុងស៊ុលហ្វាត កាបូនធាតុ សមាសធាតុលោហធាតុ និងសមាសធាតុសរីរាង្គដែលមានសម្ពាធចំហាយទឹកទាប។'
$ws.Range("G4").Value = 'PM 2.5: ពួកវាជាភាគល្អិតព្យួរដែលមានអង្កត់ផ្ចិតខ្យល់រហូតដល់ 2.5 µm ហៅថា ភាគល្អិតល្អ ឬប្រភាគល្អ។ ពួកវាជាប់ពាក់ព័ន្ធក្នុងការកើនឡើងនៃការស្លាប់ និងជំងឺដោយសារមូលហេតុផ្លូវដង្ហើម និងសរសៃឈាមបេះដូង។'
$ws.Range("H4").Value = 'PM 10: ពួកវាអាចកំណត់បានថាជាភាគល្អិតរឹង ឬរាវនៃធូលី ផេះ កំបោរ ភាគល្អិតលោហធាតុ ស៊ីម៉ងត៍ ឬលំអងដែលបែកខ្ចាត់ខ្ចាយក្នុងបរិយាកាស ហើយអង្កត់ផ្ចិតរបស់វាប្រែប្រួលចន្លោះពី 2.5 ទៅ 10 µm ។ វាអាចបណ្តាលឱ្យមានផលប៉ះពាល់ដល់ប្រព័ន្ធផ្លូវដង្ហើមរបស់មនុស្ស ប៉ុន្តែពួកវាមានះថាក់តិចជាង PM2.5 ចាប់តាំងពីមានទំហំធំជាង ពួកវាមិនអាចឆ្លងកាត់ alveoli សួតបានទេ។'
$ws.Range("I4").Value = 'កម្រិតបន្ទុកឧបករណ៍'

$ws.Range("I15").Select() | Out-Null
